$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "Mariënstede (Stichting) (incl. Vughterstede)" before row 97
# (which currently holds "Mediant, Stichting voor Geestelijke
# Gezondheidszorg Oost- en Midden Twente"), keeping the alphabetical order.
$ws.Rows("97:97").Insert()
$ws.Range("A97").Value = "Mariënstede (Stichting) (incl. Vughterstede)"
$ws.Range("B97").Value = "Vastgesteld"

# Insert "Schärwachter B.V." before what is now row 143 (originally row 142,
# holding "Severinusstichting"); the earlier insert shifted it down by one.
$ws.Rows("143:143").Insert()
$ws.Range("A143").Value = "Schärwachter B.V."
$ws.Range("B143").Value = "Vastgesteld"

# Restore the active-cell selection shown in the saved workbook.
$ws.Range("A2").Select()
